$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (tata) from "TATA 1" to "TATA 2" for rows 3 through 42
for ($row = 3; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq "TATA 1") {
        $cell.Value = "TATA 2"
    }
}
